$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "53.962.72"
$ws.Cells.Item(2, 5).Value = "  -0.66%  "

$ws.Cells.Item(3, 4).Value = "2.262.08"
$ws.Cells.Item(3, 5).Value = "  -0.59%  "

$ws.Cells.Item(4, 5).Value = "  -0.14%  "

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "495.13"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "127.35"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.06%  "

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.15%  "

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -1.47%  "

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0949"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -0.59%  "

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.152"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.33%  "

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.334"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +2.69%  "

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.47%  "

$ws.Cells.Item(13, 4).Value = "2.662.78"
$ws.Cells.Item(13, 5).Value = "  -1.44%  "

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "22.54"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.81%  "

$ws.Cells.Item(15, 4).Value = "53.944.01"
$ws.Cells.Item(15, 5).Value = "  -0.86%  "

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000129"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.32%  "

$ws.Cells.Item(17, 4).Value = "2.268.50"
$ws.Cells.Item(17, 5).Value = "  -1.65%  "

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.18"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.55%  "

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.13"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.77%  "

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "300.63"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -2.19%  "

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.28"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -3.27%  "

$ws.Cells.Item(22, 5).Value = "  +0.17%  "

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "60.76"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -3.57%  "

$ws.Cells.Item(24, 5).Value = "  +0.09%  "

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.148"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -2.06%  "

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.25"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.59%  "

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "172.68"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +1.10%  "

$ws.Cells.Item(28, 5).Value = "  -1.11%  "

$ws.Cells.Item(29, 4).Value = "0.0₃0686"
$ws.Cells.Item(29, 5).Value = "  -1.54%  "

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.90"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.51%  "

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.07"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.61%  "

$ws.Cells.Item(32, 5).Value = "  -0.04%  "

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.68"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -0.18%  "

$ws.Cells.Item(34, 5).Value = "  +0.82%  "

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.931"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +7.53%  "

$ws.Cells.Item(36, 5).Value = "  -1.77%  "

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.68"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.44%  "

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.370"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.64%  "

$ws.Cells.Item(39, 5).Value = "  -2.00%  "

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.34"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.57%  "

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "124.54"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -3.47%  "

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.77"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.27%  "

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0489"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +0.28%  "

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0886"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -1.00%  "

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.541"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.75%  "

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "236.44"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -3.89%  "

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.369"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -1.72%  "

$ws.Cells.Item(48, 5).Value = "  -0.47%  "

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "10.75"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.42%  "

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "16.06"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.19%  "

$ws.Cells.Item(51, 5).Value = "  -0.60%  "
